$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newPath = "\\fs1-cbr.nexus.csiro.au\{ev-neap}\work\extent\processing\NEAP_intermediate\NVIS_PRE1750_IUCNGET_DK_20240730.tif"
$newTarget = "file:///\\fs1-cbr.nexus.csiro.au\%7bev-neap%7d\work\extent\processing\NEAP_intermediate\NVIS_PRE1750_IUCNGET_DK_20240730.tif"

$cell = $ws.Range("B4")
$cell.Value = $newPath

$ws.Hyperlinks.Add($cell, $newTarget)
$cell.Style = "Hyperlink"

$ws.Range("A4").Select()
